$wb = $excel.ActiveWorkbook
$wsCap = $wb.Worksheets.Item("Cap")
$wsRes = $wb.Worksheets.Item("Res")

# --- Res sheet: add two new resistor rows (31 and 32) ---

# Row 31 - RES030 / Res 30k 1% 0603
$wsRes.Range("G30").Copy()
$wsRes.Range("G31").PasteSpecial(-4122)
$wsRes.Range("H30").Copy()
$wsRes.Range("H31").PasteSpecial(-4122)
$wsRes.Range("N30").Copy()
$wsRes.Range("N31").PasteSpecial(-4122)

$wsRes.Range("A31").Value = "RES030"
$wsRes.Range("D31").Value = "Res 30k 1% 0603"
$wsRes.Range("D31").ClearFormats()
$wsRes.Range("E31").Value = "Yageo"
$wsRes.Range("F31").Value = "RC0603FR-0730KL"
$wsRes.Range("G31").Value = "30k"
$wsRes.Range("H31").Value = "'1%"
$wsRes.Range("I31").Value = "1/10W"
$wsRes.Range("J31").Value = "75V"
$wsRes.Range("K31").Value = "SMT"
$wsRes.Range("L31").Value = "resistor"
$wsRes.Range("M31").Value = "TEP_243-Altium_RLC-Lib.SchLib"
$wsRes.Range("N31").Value = "R_1608[0603]"
$wsRes.Range("O31").Value = "TEP_243-Altium_RLC-Lib.PcbLib"
$wsRes.Range("P31").Value = "LCSC"
$wsRes.Range("Q31").Value = "C100001"

# Row 32 - RES031 / Res 150R 1% 0805
$wsRes.Range("G30").Copy()
$wsRes.Range("G32").PasteSpecial(-4122)
$wsRes.Range("H30").Copy()
$wsRes.Range("H32").PasteSpecial(-4122)
$wsRes.Range("N30").Copy()
$wsRes.Range("N32").PasteSpecial(-4122)

$wsRes.Range("A32").Value = "RES031"
$wsRes.Range("D32").Value = "Res 150R 1% 0805"
$wsRes.Range("D32").ClearFormats()
$wsRes.Range("E32").Value = "Yageo"
$wsRes.Range("F32").Value = "RC0805FR-07150RL"
$wsRes.Range("G32").Value = "150R"
$wsRes.Range("H32").Value = "'1%"
$wsRes.Range("I32").Value = "1/8W"
$wsRes.Range("J32").Value = "150V"
$wsRes.Range("K32").Value = "SMT"
$wsRes.Range("L32").Value = "resistor"
$wsRes.Range("M32").Value = "TEP_243-Altium_RLC-Lib.SchLib"
$wsRes.Range("N32").Value = "R_2012[0805]"
$wsRes.Range("O32").Value = "TEP_243-Altium_RLC-Lib.PcbLib"
$wsRes.Range("P32").Value = "LCSC"
$wsRes.Range("Q32").Value = "C114523"

# --- Update sheet view / selection state ---
# Cap is no longer the selected tab; scroll right to column I, select D27
$wsCap.Activate()
$wsCap.Range("D27").Select()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1

# Res becomes the active/selected tab; scroll down near row 22, select B33
$wsRes.Activate()
$wsRes.Range("B33").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
